$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Schedule table (Table 1): swap "Graduate School" (week 3) and
#    "CV/Academic resume" (week 4); rename "Open workshop" -> "Open session".
# ---------------------------------------------------------------------------
$scheduleTable = $d.Tables.Item(1)
$scheduleTable.Cell(4, 4).Range.Text = "CV/Academic resume"
$scheduleTable.Cell(5, 4).Range.Text = "Graduate School"
$scheduleTable.Cell(12, 4).Range.Text = "Open session"

# ---------------------------------------------------------------------------
# 2) Deliverables table (Table 2): update dates / product names, reorder rows
#    and add two new rows.
# ---------------------------------------------------------------------------
$delivTable = $d.Tables.Item(2)

# Row 2: Early Summer Reflection -- date moves a week later.
$delivTable.Cell(2, 2).Range.Text = "June 25th"

# Row 3: was "Research Abstract*" / "July 10th" -> becomes the
# "Updated CV/Academic Resume" / "July 9th" deliverable.
$delivTable.Cell(3, 1).Range.Text = "Updated CV/Academic Resume"
$delivTable.Cell(3, 2).Range.Text = "July 9th"

# Row 4: Midsummer Reflection -- date moves later.
$delivTable.Cell(4, 2).Range.Text = "July 23rd"

# Rows 5 & 6: both "August 1st (or by arrangement)" deadlines slip to the 6th.
$delivTable.Cell(5, 2).Range.Text = "August 6th (or by arrangement)"
$delivTable.Cell(6, 2).Range.Text = "August 6th (or by arrangement)"

# Row 7: was "Draft Presentation" / "August 17th" -> becomes the new
# "Coding project" / "August 13th" deliverable.
$delivTable.Cell(7, 1).Range.Text = "Coding project"
$delivTable.Cell(7, 2).Range.Text = "August 13th"

# Row 8: was "Final Presentation" / "August 21th" -> becomes
# "Draft Presentation" / "August 17th" (shifted down from row 7).
$delivTable.Cell(8, 1).Range.Text = "Draft Presentation"
$delivTable.Cell(8, 2).Range.Text = "August 17th"

# Insert a new row before the current row 9 ("End of Summer Reflection") to
# hold "Final Presentation" / "August 20th" (shifted down from row 8).
$refRow = $delivTable.Rows.Item(9)
$newRow1 = $delivTable.Rows.Add($refRow)
$newRow1.Cells.Item(1).Split(1, 2)
$newRow1.Cells.Item(1).Range.Text = "Final Presentation"
$newRow1.Cells.Item(1).Range.Paragraphs.Item(1).Style = "Compact"
$newRow1.Cells.Item(2).Range.Text = "August 20th"
$newRow1.Cells.Item(2).Range.Paragraphs.Item(1).Style = "Compact"

# Row 10 is now "End of Summer Reflection" / "August 27th" (unchanged).

# Append a new row at the end of the table for
# "Research Abstract*" / "September 11th (ABRCMS site)".
$newRow2 = $delivTable.Rows.Add()
$newRow2.Cells.Item(1).Split(1, 2)
$newRow2.Cells.Item(1).Range.Text = "Research Abstract*"
$newRow2.Cells.Item(1).Range.Paragraphs.Item(1).Style = "Compact"
$newRow2.Cells.Item(2).Range.Text = "September 11th (ABRCMS site)"
$newRow2.Cells.Item(2).Range.Paragraphs.Item(1).Style = "Compact"

# ---------------------------------------------------------------------------
# 3) Body text tweak: "deliverable" -> "deliverables".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Please email deliverable to", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Please email deliverables to", 2)
